$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 474.2
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1226

$ws.Range("H64").Value = 3883.7307
$ws.Range("I64").Value = 3515.389
$ws.Range("J64").Value = 4712.5
$ws.Range("K64").Value = 3515.389
$ws.Range("L64").Value = 4712.5
$ws.Range("M64").Value = -3267.389
$ws.Range("N64").Value = -5208.5

$ws.Range("H67").Value = 3883.7307
$ws.Range("I67").Value = 3515.389
$ws.Range("J67").Value = 4712.5
$ws.Range("K67").Value = 3515.389
$ws.Range("L67").Value = 4712.5
$ws.Range("M67").Value = -2657.389
$ws.Range("N67").Value = -6428.5

$ws.Range("H137").Value = 2754.6863
$ws.Range("I137").Value = 1565.3043
$ws.Range("J137").Value = 3731.6785
$ws.Range("K137").Value = 4695.9129
$ws.Range("L137").Value = 11195.0355
$ws.Range("M137").Value = -2145.9129
$ws.Range("N137").Value = -16295.0355

$ws.Range("H138").Value = 1342023
$ws.Range("I138").Value = 3858.9092
$ws.Range("J138").Value = 1600265.2
$ws.Range("K138").Value = 11576.7276
$ws.Range("L138").Value = 4800795.6
$ws.Range("M138").Value = -6436.7276
$ws.Range("N138").Value = -4811075.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 53339.668
$ws.Range("J23").Value = 60007
$ws.Range("L23").Value = 60007
$ws.Range("N23").Value = -60525

$ws.Range("H45").Value = 1591.8928
$ws.Range("I45").Value = 1587.3182
$ws.Range("J45").Value = 1608.6666
$ws.Range("K45").Value = 1587.3182
$ws.Range("L45").Value = 1608.6666
$ws.Range("M45").Value = -1210.3182
$ws.Range("N45").Value = -2362.6666

$ws.Range("H74").Value = 3715.375
$ws.Range("I74").Value = 1779.4375
$ws.Range("K74").Value = 1779.4375
$ws.Range("M74").Value = -905.4375

$ws.Range("H77").Value = 3715.375
$ws.Range("I77").Value = 1779.4375
$ws.Range("K77").Value = 8897.1875
$ws.Range("M77").Value = -4529.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22483.256
$ws.Range("I134").Value = 2648.3784
$ws.Range("K134").Value = 7945.135200000001
$ws.Range("M134").Value = -5410.135200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9000
$ws.Range("J25").Value = 9500
$ws.Range("L25").Value = 9500
$ws.Range("N25").Value = -9848

$ws.Range("H31").Value = 4992.522
$ws.Range("I31").Value = 2090
$ws.Range("K31").Value = 2090
$ws.Range("M31").Value = -1795

$ws.Range("H34").Value = 4992.522
$ws.Range("I34").Value = 2090
$ws.Range("K34").Value = 2090
$ws.Range("M34").Value = -1888

$ws.Range("H62").Value = 2900.4443
$ws.Range("I62").Value = 2872
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2872
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2248
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 2900.4443
$ws.Range("I65").Value = 2872
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14360
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11240
$ws.Range("N65").Value = -21240

$ws.Range("H99").Value = 2612.5
$ws.Range("I99").Value = 2650
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2650
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -1152
$ws.Range("N99").Value = -5496

$ws.Range("H126").Value = 2612.5
$ws.Range("I126").Value = 2650
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 7950
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -5480
$ws.Range("N126").Value = -12440

$ws.Range("H132").Value = 4008.3572
$ws.Range("I132").Value = 4288.125
$ws.Range("K132").Value = 12864.375
$ws.Range("M132").Value = -10334.375

$ws.Range("H134").Value = 2877.25
$ws.Range("I134").Value = 1791.9706
$ws.Range("J134").Value = 3962.5293
$ws.Range("K134").Value = 5375.9118
$ws.Range("L134").Value = 11887.5879
$ws.Range("M134").Value = -2840.9118
$ws.Range("N134").Value = -16957.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 846.3333
$ws.Range("I114").Value = 769.5
$ws.Range("K114").Value = 2308.5
$ws.Range("M114").Value = 945.5

$ws.Range("H117").Value = 1090.8
$ws.Range("J117").Value = 1090.8
$ws.Range("L117").Value = 3272.4
$ws.Range("N117").Value = -10156.4

$ws.Range("H122").Value = 1294.0238
$ws.Range("J122").Value = 1670.2222
$ws.Range("L122").Value = 15031.9998
$ws.Range("N122").Value = -19931.9998

$ws.Range("H129").Value = 2069.318
$ws.Range("J129").Value = 1411.7858
$ws.Range("L129").Value = 4235.357400000001
$ws.Range("N129").Value = -14235.3574

$ws.Range("H131").Value = 545.24243
$ws.Range("I131").Value = 282.45764
$ws.Range("J131").Value = 932.85
$ws.Range("K131").Value = 847.37292
$ws.Range("L131").Value = 2798.55
$ws.Range("M131").Value = 4192.62708
$ws.Range("N131").Value = -12878.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6211.1665
$ws.Range("J102").Value = 6717.2856
$ws.Range("L102").Value = 6717.2856
$ws.Range("N102").Value = -9961.285599999999

$ws.Range("H122").Value = 11724.875
$ws.Range("I122").Value = 11160
$ws.Range("J122").Value = 12666.333
$ws.Range("K122").Value = 33480
$ws.Range("L122").Value = 37998.999
$ws.Range("M122").Value = -31030
$ws.Range("N122").Value = -42898.999

$ws.Range("H126").Value = 2994.0476
$ws.Range("I126").Value = 1925
$ws.Range("K126").Value = 5775
$ws.Range("M126").Value = -3305

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4737.5
$ws.Range("I7").Value = 2883.9473
$ws.Range("K7").Value = 2883.9473
$ws.Range("M7").Value = -2771.9473

$ws.Range("H40").Value = 4019.95
$ws.Range("I40").Value = 3729.353
$ws.Range("J40").Value = 5666.6665
$ws.Range("K40").Value = 3729.353
$ws.Range("L40").Value = 5666.6665
$ws.Range("M40").Value = -3593.353
$ws.Range("N40").Value = -5938.6665

$ws.Range("H55").Value = 444561.44
$ws.Range("I55").Value = 666816.7
$ws.Range("J55").Value = 51
$ws.Range("K55").Value = 666816.7
$ws.Range("L55").Value = 51
$ws.Range("M55").Value = -666643.7
$ws.Range("N55").Value = -397

$ws.Range("H82").Value = 1916.6
$ws.Range("J82").Value = 3302.8572
$ws.Range("L82").Value = 3302.8572
$ws.Range("N82").Value = -4024.8572

$ws.Range("H85").Value = 1916.6
$ws.Range("J85").Value = 3302.8572
$ws.Range("L85").Value = 3302.8572
$ws.Range("N85").Value = -5798.8572

$ws.Range("H122").Value = 5281.7812
$ws.Range("I122").Value = 3987.261
$ws.Range("J122").Value = 8590
$ws.Range("K122").Value = 11961.783
$ws.Range("L122").Value = 25770
$ws.Range("M122").Value = -9511.782999999999
$ws.Range("N122").Value = -30670

$ws.Range("H126").Value = 4737.5
$ws.Range("I126").Value = 2883.9473
$ws.Range("K126").Value = 8651.841899999999
$ws.Range("M126").Value = -6181.841899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1720.8684
$ws.Range("I132").Value = 1189.16
$ws.Range("J132").Value = 2743.3845
$ws.Range("K132").Value = 3567.48
$ws.Range("L132").Value = 8230.1535
$ws.Range("M132").Value = -1037.48
$ws.Range("N132").Value = -13290.1535

$ws.Range("H136").Value = 4375.21
$ws.Range("I136").Value = 2250.653
$ws.Range("K136").Value = 6751.958999999999
$ws.Range("M136").Value = -4201.958999999999
